# Generate Report for Handback
# Marks the two handed-off files as handed back: updates the status text,
# fills in the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns (with hyperlinks on the target file column)
# on both the zh-cn and de-de language sheets, and widens the columns that
# now hold longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$fileA = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$fileB = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"

$fileA_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff208caef7e72897025f5f3223733297529e9304/e2e/2fba2052-50bf-4790-8c8c-2e07826b5dc7.md"
$fileB_url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff208caef7e72897025f5f3223733297529e9304/e2e/68e22ea2-5fb0-443a-9dfc-bf9dd7778106.md"

$zhcn_handbackA = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.ef624036b7a2d9357573e67962c90ab2b542cc3e.zh-cn.xlf"
$zhcn_handbackB = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.1e28be1cea68a5b58d1cc99c5bc006ea64757a9e.zh-cn.xlf"
$dede_handbackA = "2fba2052-50bf-4790-8c8c-2e07826b5dc7.ef624036b7a2d9357573e67962c90ab2b542cc3e.de-de.xlf"
$dede_handbackB = "68e22ea2-5fb0-443a-9dfc-bf9dd7778106.1e28be1cea68a5b58d1cc99c5bc006ea64757a9e.de-de.xlf"

$zhcn_handbackDate = "2016-08-24 03:01:17"
$dede_handbackDate = "2016-08-24 03:01:24"

$hyperlinkColor = 15570276  # RGB(100,149,237) == the workbook's existing "HyperLink" style color (FF6495ED)

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# 1. Status: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3 and the "Status" column on both language sheets all
#    share this text.)
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows, and hyperlink the target
#    file cells.
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = $fileA
$zhcn.Range("J2").Value = $zhcn_handbackA
$zhcn.Range("K2").Value = $zhcn_handbackDate

$zhcn.Range("I3").Value = $fileB
$zhcn.Range("J3").Value = $zhcn_handbackB
$zhcn.Range("K3").Value = $zhcn_handbackDate

# ---------------------------------------------------------------------
# 3. de-de sheet: same as above.
# ---------------------------------------------------------------------
$dede.Range("I2").Value = $fileA
$dede.Range("J2").Value = $dede_handbackA
$dede.Range("K2").Value = $dede_handbackDate

$dede.Range("I3").Value = $fileB
$dede.Range("J3").Value = $dede_handbackB
$dede.Range("K3").Value = $dede_handbackDate

# ---------------------------------------------------------------------
# 4. Re-create the hyperlinks on both sheets so that A2, I2, A3, I3 end
#    up in that order (I2/I3 link to the same .md files as A2/A3).
# ---------------------------------------------------------------------
foreach ($ws in @($zhcn, $dede)) {
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $fileA_url, "", "", $fileA) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), $fileA_url, "", "", $fileA) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $fileB_url, "", "", $fileB) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $fileB_url, "", "", $fileB) | Out-Null

    Style-AsHyperlink $ws.Range("I2")
    Style-AsHyperlink $ws.Range("I3")
}

# ---------------------------------------------------------------------
# 5. Column widths: widen the (now longer-content) Status / Latest Target
#    File / Latest Handback File columns.
# ---------------------------------------------------------------------
$overview.Range("E1").ColumnWidth = 29.084
$overview.Range("F1").ColumnWidth = 29.084

foreach ($ws in @($zhcn, $dede)) {
    $ws.Range("C1").ColumnWidth = 29.084
    $ws.Range("I1").ColumnWidth = 39.17
    $ws.Range("J1").ColumnWidth = 39.17
}
